$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text formatting (matches original inlineStr text cells)
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "E50", "D51", "E51")
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.294.19'
$ws.Range("E2").Value = '  -2.50%  '
$ws.Range("D3").Value = '2.970.12'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '588.66'
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").Value = '141.79'
$ws.Range("E6").Value = '  -5.49%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("D9").Value = '2.969.50'
$ws.Range("E9").Value = '  -2.39%  '
$ws.Range("E10").Value = '  -5.81%  '
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  +2.25%  '
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  -3.06%  '
$ws.Range("D14").Value = '33.90'
$ws.Range("E14").Value = '  -5.13%  '
$ws.Range("D16").Value = '3.458.83'
$ws.Range("E16").Value = '  -2.38%  '
$ws.Range("D17").Value = '6.99'
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").Value = '61.267.55'
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("D19").Value = '2.965.68'
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("D20").Value = '447.73'
$ws.Range("E20").Value = '  -5.76%  '
$ws.Range("D21").Value = '13.84'
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").Value = '7.31'
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("D24").Value = '81.10'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '12.08'
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("D26").Value = '2.16'
$ws.Range("E26").Value = '  -8.19%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '9.96'
$ws.Range("E28").Value = '  -4.60%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '2.63'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '6.84'
$ws.Range("E31").Value = '  -5.79%  '
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -5.63%  '
$ws.Range("D33").Value = '26.84'
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("E35").Value = '  -4.62%  '
$ws.Range("D36").Value = '0.0₃0774'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").Value = '5.69'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").Value = '50.08'
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("D40").Value = '9.12'
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").Value = '0.119'
$ws.Range("E41").Value = '  +5.01%  '
$ws.Range("D42").Value = '2.76'
$ws.Range("E42").Value = '  -10.08%  '
$ws.Range("D43").Value = '386.27'
$ws.Range("E43").Value = '  -8.59%  '
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").Value = '2.697.68'
$ws.Range("E45").Value = '  -4.43%  '
$ws.Range("D46").Value = '0.263'
$ws.Range("E46").Value = '  -7.26%  '
$ws.Range("D47").Value = '36.97'
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("D48").Value = '130.37'
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '2.15'
$ws.Range("E51").Value = '  -0.65%  '
